$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'57.700.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  -0.61%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'2.413.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  -1.63%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'510.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -2.82%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'133.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +1.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -0.56%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Formula = "'0.558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'  -1.52%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'2.452.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  -0.21%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'0.0981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +0.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Formula = "'  -1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Formula = "'  -0.94%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'4.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Formula = "'2.849.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  -1.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'57.553.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  -0.76%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'21.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +0.66%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +0.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'2.447.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  -0.06%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'10.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  -1.57%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'4.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +0.24%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'315.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  +0.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Formula = "'  +4.83%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -0.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'5.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  -2.06%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'65.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  +0.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'0.994"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  -0.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'2.549.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  -0.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Formula = "'  -0.84%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'0.382"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  -5.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Formula = "'7.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +4.29%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'173.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Formula = "'0.0₃0736"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  -0.39%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Formula = "'  -0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Formula = "'6.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Formula = "'1.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  -0.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  -0.18%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'0.991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  -0.55%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'18.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  +0.87%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Formula = "'1.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  +4.88%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'3.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +1.17%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'36.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  +1.26%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'0.814"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -0.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'1.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +1.32%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'135.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +10.93%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -0.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'5.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'  +4.42%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'258.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  -1.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'0.574"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -2.02%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.0919"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  -0.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'0.0493"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  -0.53%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'0.0215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +1.28%  "
$ws.Range("E51").Style = "Normal"
